$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 5960.4
$ws.Range("J19").Value = 5960.4
$ws.Range("L19").Value = 5960.4
$ws.Range("N19").Value = -6310.4

$ws.Range("H40").Value = 2916.5
$ws.Range("I40").Value = 3099.8
$ws.Range("J40").Value = 2000
$ws.Range("K40").Value = 3099.8
$ws.Range("L40").Value = 2000
$ws.Range("M40").Value = -2924.8
$ws.Range("N40").Value = -2350

$ws.Range("H43").Value = 412889.3
$ws.Range("I43").Value = 4949
$ws.Range("K43").Value = 4949
$ws.Range("M43").Value = -4880

$ws.Range("H107").Value = 23438376
$ws.Range("I107").Value = 8929002
$ws.Range("J107").Value = 125004000
$ws.Range("K107").Value = 8929002
$ws.Range("L107").Value = 125004000
$ws.Range("M107").Value = -8927082
$ws.Range("N107").Value = -125007840

$ws.Range("H123").Value = 91400.164
$ws.Range("J123").Value = 91400.164
$ws.Range("L123").Value = 91400.164
$ws.Range("N123").Value = -101200.164

$ws.Range("H132").Value = 1548.8113
$ws.Range("I132").Value = 1274.0952
$ws.Range("J132").Value = 2597.7273
$ws.Range("K132").Value = 3822.2856
$ws.Range("L132").Value = 7793.1819
$ws.Range("M132").Value = -1292.2856
$ws.Range("N132").Value = -12853.1819

$ws.Range("H138").Value = 5459.5137
$ws.Range("J138").Value = 7263.0835
$ws.Range("L138").Value = 21789.2505
$ws.Range("N138").Value = -32069.2505

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 50
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 50
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 50
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -274

$ws.Range("H28").Value = 17603.25
$ws.Range("I28").Value = 5490.3335
$ws.Range("K28").Value = 5490.3335
$ws.Range("M28").Value = -5298.3335

$ws.Range("H45").Value = 6132
$ws.Range("I45").Value = 2122.7
$ws.Range("J45").Value = 12814.167
$ws.Range("K45").Value = 2122.7
$ws.Range("L45").Value = 12814.167
$ws.Range("M45").Value = -1745.7
$ws.Range("N45").Value = -13568.167

$ws.Range("H99").Value = 17603.25
$ws.Range("I99").Value = 5490.3335
$ws.Range("K99").Value = 5490.3335
$ws.Range("M99").Value = -2495.3335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 50
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 50
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -280

$ws.Range("H22").Value = 329.7
$ws.Range("I22").Value = 299.7143
$ws.Range("J22").Value = 399.66666
$ws.Range("K22").Value = 299.7143
$ws.Range("L22").Value = 399.66666
$ws.Range("M22").Value = -126.7143
$ws.Range("N22").Value = -745.66666

$ws.Range("H134").Value = 4282.6094
$ws.Range("I134").Value = 1620.0731
$ws.Range("K134").Value = 4860.219300000001
$ws.Range("M134").Value = -2325.219300000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4558.9
$ws.Range("I16").Value = 2790.0715
$ws.Range("K16").Value = 2790.0715
$ws.Range("M16").Value = -2503.0715

$ws.Range("H58").Value = 10644535
$ws.Range("I58").Value = 22729962
$ws.Range("K58").Value = 22729962
$ws.Range("M58").Value = -22729759

$ws.Range("H99").Value = 8316
$ws.Range("I99").Value = 8468.799999999999
$ws.Range("K99").Value = 8468.799999999999
$ws.Range("M99").Value = -6970.799999999999

$ws.Range("H113").Value = 4558.9
$ws.Range("I113").Value = 2790.0715
$ws.Range("K113").Value = 2790.0715
$ws.Range("M113").Value = -620.0715

$ws.Range("H122").Value = 1892.2858
$ws.Range("I122").Value = 1849.2
$ws.Range("K122").Value = 5547.6
$ws.Range("M122").Value = -3097.6

$ws.Range("H126").Value = 8316
$ws.Range("I126").Value = 8468.799999999999
$ws.Range("K126").Value = 25406.4
$ws.Range("M126").Value = -22936.4

$ws.Range("H132").Value = 5267.027
$ws.Range("I132").Value = 1961.2858
$ws.Range("J132").Value = 9605.8125
$ws.Range("K132").Value = 5883.857400000001
$ws.Range("L132").Value = 28817.4375
$ws.Range("M132").Value = -3353.857400000001
$ws.Range("N132").Value = -33877.4375

$ws.Range("H134").Value = 6670.5
$ws.Range("I134").Value = 5659.8047
$ws.Range("K134").Value = 16979.4141
$ws.Range("M134").Value = -14444.4141

$ws.Range("H136").Value = 10644535
$ws.Range("I136").Value = 22729962
$ws.Range("K136").Value = 68189886
$ws.Range("M136").Value = -68187336

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 8000
$ws.Range("J57").Value = 8000
$ws.Range("L57").Value = 24000
$ws.Range("N57").Value = -25118

$ws.Range("H131").Value = 2026.6897
$ws.Range("I131").Value = 1225.2
$ws.Range("J131").Value = 2193.6667
$ws.Range("K131").Value = 3675.6
$ws.Range("L131").Value = 6581.000100000001
$ws.Range("M131").Value = 1364.4
$ws.Range("N131").Value = -16661.0001

$ws.Range("H140").Value = 223349.89
$ws.Range("I140").Value = 223349.89
$ws.Range("K140").Value = 670049.67
$ws.Range("M140").Value = -664869.67

$ws.Range("H141").Value = 8757.166999999999
$ws.Range("I141").Value = 3512.2856
$ws.Range("J141").Value = 16100
$ws.Range("K141").Value = 10536.8568
$ws.Range("L141").Value = 48300
$ws.Range("M141").Value = -5356.856800000001
$ws.Range("N141").Value = -58660

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4267.3335
$ws.Range("I126").Value = 4267.3335
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 12802.0005
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -10332.0005
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 10571.857
$ws.Range("I93").Value = 10800.6
$ws.Range("J93").Value = 10000
$ws.Range("K93").Value = 10800.6
$ws.Range("L93").Value = 10000
$ws.Range("M93").Value = -9552.6
$ws.Range("N93").Value = -12496

$ws.Range("H132").Value = 10007202
$ws.Range("I132").Value = 16671151
$ws.Range("K132").Value = 50013453
$ws.Range("M132").Value = -50010923

$ws.Range("H136").Value = 9504.700000000001
$ws.Range("I136").Value = 2239.111
$ws.Range("K136").Value = 6717.333
$ws.Range("M136").Value = -4167.333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 254165.94
$ws.Range("I122").Value = 336287.9
$ws.Range("K122").Value = 1008863.7
$ws.Range("M122").Value = -1006413.7

$ws.Range("H132").Value = 8780476
$ws.Range("I132").Value = 12198977
$ws.Range("K132").Value = 36596931
$ws.Range("M132").Value = -36594401

$ws.Range("H136").Value = 19630834
$ws.Range("I136").Value = 38462236
$ws.Range("K136").Value = 115386708
$ws.Range("M136").Value = -115384158
